$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - account holder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number and surname
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 31.01.2025"

# Row 6 - transaction 1
$ws.Range("B6").Value = "01.02."
$ws.Range("C6").Value = "02.02."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-31191235"
$ws.Range("E6").Value = "55,63-"

# Row 7 - transaction 2
$ws.Range("B7").Value = "04.02."
$ws.Range("C7").Value = "05.02."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 32748845"
$ws.Range("E7").Value = "40,49-"

# Row 8 - transaction 3
$ws.Range("B8").Value = "06.02."
$ws.Range("C8").Value = "07.02."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 38500995"
$ws.Range("E8").Value = "87,53-"

# Row 9 - transaction 4
$ws.Range("B9").Value = "08.02."
$ws.Range("C9").Value = "09.02."
$ws.Range("D9").Value = "PAYPAL CWEPLB"
$ws.Range("E9").Value = "85,38-"

# Row 10 - transaction 5
$ws.Range("B10").Value = "10.02."
$ws.Range("C10").Value = "11.02."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "25,39-"

# Row 11 - now empty (one fewer transaction), E11 style changes to match the
# "right + vertical-center + wrap" alignment used elsewhere (style 12).
# (E11 is part of a merged cell E11:F11, so ClearContents doesn't take
# effect on it reliably - assign an empty string value instead.)
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Row 12 - closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 13.02.2025"
$ws.Range("E12").Value = "294,42-"

# Row 13 - next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 19.02.2025"
